$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    ,@(23, 0.7289727926254272, 128, 14500, 1, 0.5990680456161499, 0.001, "Wed Jul 22 17:35:12 2020")
    ,@(24, 0.9688101410865784, 128, 14500, 25, 0.5233844518661499, 0.0001, "Wed Jul 22 17:37:45 2020")
    ,@(25, 0.9695556163787842, 128, 14500, 30, 0.4996465444564819, 0.0001, "Wed Jul 22 17:55:20 2020")
    ,@(26, 0.9646815657615662, 128, 14500, 30, 0.5375465154647827, 0.0001, "Wed Jul 22 18:02:30 2020")
    ,@(27, 0.9731359481811523, 128, 14500, 30, 0.4832836389541626, 0.0001, "Wed Jul 22 18:58:16 2020")
    ,@(28, 0.9776589870452881, 128, 14500, 30, 0.4949248135089874, 0.0001, "Wed Jul 22 18:59:43 2020")
    ,@(29, 0.9766995906829834, 128, 14500, 30, 0.5067653656005859, 0.0001, "Wed Jul 22 19:25:35 2020")
    ,@(30, 0.7157986164093018, 128, 14500, 1, 0.6001615524291992, 0.0001, "Thu Jul 23 13:56:45 2020")
    ,@(31, 0.7362847328186035, 128, 14500, 1, 0.6436253786087036, 0.0001, "Thu Jul 23 14:07:07 2020")
    ,@(32, 0.8692708611488342, 128, 14500, 3, 0.5539892911911011, 0.0001, "Thu Jul 23 14:07:32 2020")
    ,@(33, 0.8687499761581421, 128, 14500, 3, 0.575569212436676, 0.0001, "Thu Jul 23 15:08:05 2020")
    ,@(34, 0.8763889074325562, 128, 14500, 3, 0.5662683844566345, 0.0001, "Thu Jul 23 15:12:19 2020")
    ,@(35, 0.8685764074325562, 128, 14500, 3, 0.5539690256118774, 0.0001, "Thu Jul 23 15:13:24 2020")
    ,@(36, 0.8607639074325562, 128, 14500, 3, 0.5884552001953125, 0.0001, "Thu Jul 23 15:14:46 2020")
    ,@(37, 0.9710069298744202, 128, 14500, 25, 0.5304288864135742, 0.0001, "Thu Jul 23 15:17:39 2020")
    ,@(38, 0.7401041388511658, 128, 14500, 1, 0.610766589641571, 0.0001, "Thu Jul 23 17:00:02 2020")
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
    $ws.Cells.Item($r, 5).Value = $row[5]
    $ws.Cells.Item($r, 6).Value = $row[6]
    $ws.Cells.Item($r, 7).Value = $row[7]
}

Write-Host "Added rows 23-38"
